$d = $word.ActiveDocument

# --- Step 1: change "What are the constraints?" to the new constraints text ---
$d.Content.Find.Execute(
    "What are the constraints?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Constraints are that we must go back to one every time.", 2) | Out-Null

# --- Step 2: change "What are the sub-goals?" to the new "Find the finger" text ---
$d.Content.Find.Execute(
    "What are the sub-goals?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Find the finger the easiest way.  ", 2) | Out-Null

# --- Step 3: move the _GoBack bookmark from the end of the "...each time." ---
# --- paragraph to the end of the new "Find the finger..." paragraph.      ---

# Locate the paragraph that now holds "Find the finger the easiest way..."
$count = $d.Paragraphs.Count
$targetParaIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Find the finger the easiest way.*") {
        $targetParaIdx = $i
        break
    }
}
if ($targetParaIdx -eq -1) {
    throw "Could not locate the 'Find the finger the easiest way' paragraph"
}

$p = $d.Paragraphs.Item($targetParaIdx)

# Insert a throwaway character right before the paragraph mark so that the
# collapsed insertion point we need (immediately after the run's text) is no
# longer the literal "end-of-paragraph" offset.
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$insertPoint.InsertAfter("Z")

# Recompute the paragraph range (it grew by one character) and add/relocate
# the _GoBack bookmark right before that throwaway character. Re-adding a
# bookmark named "_GoBack" moves the existing one rather than duplicating it.
$p = $d.Paragraphs.Item($targetParaIdx)
$bmPos = $p.Range.End - 2
$bmTarget = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmTarget) | Out-Null

# Remove the throwaway character again; the bookmark stays anchored in place.
$dummyRange = $d.Range($bmPos, $bmPos + 1)
$dummyRange.Delete()
